# 🔄 MAJ automatique BRVM via GitHub Actions
# Refresh the "Recommandations" and "Top_YTD" sheets with the latest BRVM
# market data: updated variation/YTD figures, two new "(**)" total-return
# index rows inserted among the indices block, several stock rows
# re-ordered/re-valued, and four additional stock rows (NESTLE, FILTISAC,
# SAFCA, SITAB) appended at the bottom of the recommendations table.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Recommandations ---
$ws1 = $wb.Worksheets.Item("Recommandations")

# Row 2: BRVM - CONSOMMATION DISCRETIONNAIRE
$ws1.Range("A2").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Range("B2").Value = 0
$ws1.Range("C2").Value = 4
$ws1.Range("D2").Value = 718.56
$ws1.Range("E2").Value = 178.71
$ws1.Range("F2").Value = "🟡 Observer"
$ws1.Range("G2").Value = "➖ Neutre"

# Row 3: BRVM - INDUSTRIELS
$ws1.Range("A3").Value = "BRVM - INDUSTRIELS"
$ws1.Range("B3").Value = 0
$ws1.Range("C3").Value = 4
$ws1.Range("D3").Value = 633.87
$ws1.Range("E3").Value = 161.02
$ws1.Range("F3").Value = "🟡 Observer"
$ws1.Range("G3").Value = "➖ Neutre"

# Row 4: BRVM - SERVICES FINANCIERS
$ws1.Range("A4").Value = "BRVM - SERVICES FINANCIERS"
$ws1.Range("B4").Value = 0
$ws1.Range("C4").Value = 4
$ws1.Range("D4").Value = 622.26
$ws1.Range("E4").Value = 156.69
$ws1.Range("F4").Value = "🟡 Observer"
$ws1.Range("G4").Value = "➖ Neutre"

# Row 5: BRVM-PRESTIGE
$ws1.Range("A5").Value = "BRVM-PRESTIGE"
$ws1.Range("B5").Value = 0
$ws1.Range("C5").Value = 4
$ws1.Range("D5").Value = 596.38
$ws1.Range("E5").Value = 150.34
$ws1.Range("F5").Value = "🟡 Observer"
$ws1.Range("G5").Value = "➖ Neutre"

# Row 6: BRVM - SERVICES PUBLICS
$ws1.Range("A6").Value = "BRVM - SERVICES PUBLICS"
$ws1.Range("B6").Value = 0
$ws1.Range("C6").Value = 4
$ws1.Range("D6").Value = 471.92
$ws1.Range("E6").Value = 118.13
$ws1.Range("F6").Value = "🟡 Observer"
$ws1.Range("G6").Value = "➖ Neutre"

# Row 7: BRVM - ENERGIE
$ws1.Range("A7").Value = "BRVM - ENERGIE"
$ws1.Range("B7").Value = 0
$ws1.Range("C7").Value = 4
$ws1.Range("D7").Value = 471.42
$ws1.Range("E7").Value = 118.38
$ws1.Range("F7").Value = "🟡 Observer"
$ws1.Range("G7").Value = "➖ Neutre"

# Row 8: BRVM - TELECOMMUNICATIONS
$ws1.Range("A8").Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Range("B8").Value = 0
$ws1.Range("C8").Value = 4
$ws1.Range("D8").Value = 384.65
$ws1.Range("E8").Value = 97.32
$ws1.Range("F8").Value = "🟡 Observer"
$ws1.Range("G8").Value = "➖ Neutre"

# Row 9: BRVM-PRINCIPAL     (**)
$ws1.Range("A9").Value = "BRVM-PRINCIPAL     (**)"
$ws1.Range("B9").Value = 0
$ws1.Range("C9").Value = 1
$ws1.Range("D9").Value = 236.54
$ws1.Range("E9").Value = 236.54
$ws1.Range("F9").Value = "🟡 Observer"
$ws1.Range("G9").Value = "➖ Neutre"

# Row 10: BRVM-PRINCIPAL
$ws1.Range("A10").Value = "BRVM-PRINCIPAL"
$ws1.Range("B10").Value = 0
$ws1.Range("C10").Value = 1
$ws1.Range("D10").Value = 236.03
$ws1.Range("E10").Value = 236.03
$ws1.Range("F10").Value = "🟡 Observer"
$ws1.Range("G10").Value = "➖ Neutre"

# Row 11: BRVM - CONSOMMATION DE BASE
$ws1.Range("A11").Value = "BRVM - CONSOMMATION DE BASE"
$ws1.Range("B11").Value = 0
$ws1.Range("C11").Value = 1
$ws1.Range("D11").Value = 233.7
$ws1.Range("E11").Value = 233.7
$ws1.Range("F11").Value = "🟡 Observer"
$ws1.Range("G11").Value = "➖ Neutre"

# Row 12: BRVM - CONSOMMATION DE BASE     (**)
$ws1.Range("A12").Value = "BRVM - CONSOMMATION DE BASE     (**)"
$ws1.Range("B12").Value = 0
$ws1.Range("C12").Value = 1
$ws1.Range("D12").Value = 230.85
$ws1.Range("E12").Value = 230.85
$ws1.Range("F12").Value = "🟡 Observer"
$ws1.Range("G12").Value = "➖ Neutre"

# Row 13: BRVM – COMPOSITE TOTAL RETURN     (**)
$ws1.Range("A13").Value = "BRVM – COMPOSITE TOTAL RETURN     (**)"
$ws1.Range("B13").Value = 0
$ws1.Range("C13").Value = 1
$ws1.Range("D13").Value = 140.62
$ws1.Range("E13").Value = 140.62
$ws1.Range("F13").Value = "🟡 Observer"
$ws1.Range("G13").Value = "➖ Neutre"

# Row 14: BRVM – COMPOSITE TOTAL RETURN
$ws1.Range("A14").Value = "BRVM – COMPOSITE TOTAL RETURN"
$ws1.Range("B14").Value = 0
$ws1.Range("C14").Value = 1
$ws1.Range("D14").Value = 139.43
$ws1.Range("E14").Value = 139.43
$ws1.Range("F14").Value = "🟡 Observer"
$ws1.Range("G14").Value = "➖ Neutre"

# Row 15: EVIOSYS PACKAGING SIEM CI (SEMC)
$ws1.Range("A15").Value = "EVIOSYS PACKAGING SIEM CI (SEMC)"
$ws1.Range("B15").Value = 3
$ws1.Range("C15").Value = 0
$ws1.Range("D15").Value = 21.91
$ws1.Range("E15").Value = 7.23
$ws1.Range("F15").Value = "🟢 Achat"
$ws1.Range("G15").Value = "✅ Renforcer"

# Row 16: ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)
$ws1.Range("A16").Value = "ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)"
$ws1.Range("B16").Value = 3
$ws1.Range("C16").Value = 1
$ws1.Range("D16").Value = 13.67
$ws1.Range("E16").Value = -4.65
$ws1.Range("F16").Value = "🟢 Achat"
$ws1.Range("G16").Value = "✅ Renforcer"

# Row 17: SMB CI (SMBC)
$ws1.Range("A17").Value = "SMB CI (SMBC)"
$ws1.Range("B17").Value = 2
$ws1.Range("C17").Value = 0
$ws1.Range("D17").Value = 10.86
$ws1.Range("E17").Value = 7.41
$ws1.Range("F17").Value = "🟡 Observer"
$ws1.Range("G17").Value = "➖ Neutre"

# Row 18: UNILEVER CI (UNLC)
$ws1.Range("A18").Value = "UNILEVER CI (UNLC)"
$ws1.Range("B18").Value = 1
$ws1.Range("C18").Value = 0
$ws1.Range("D18").Value = 6.36
$ws1.Range("E18").Value = 6.36
$ws1.Range("F18").Value = "🟡 Observer"
$ws1.Range("G18").Value = "➖ Neutre"

# Row 19: NEI-CEDA CI (NEIC)
$ws1.Range("A19").Value = "NEI-CEDA CI (NEIC)"
$ws1.Range("B19").Value = 1
$ws1.Range("C19").Value = 0
$ws1.Range("D19").Value = 5.81
$ws1.Range("E19").Value = 5.81
$ws1.Range("F19").Value = "🟡 Observer"
$ws1.Range("G19").Value = "➖ Neutre"

# Row 20: SERVAIR ABIDJAN CI (ABJC)
$ws1.Range("A20").Value = "SERVAIR ABIDJAN CI (ABJC)"
$ws1.Range("B20").Value = 1
$ws1.Range("C20").Value = 0
$ws1.Range("D20").Value = 5.17
$ws1.Range("E20").Value = 5.17
$ws1.Range("F20").Value = "🟡 Observer"
$ws1.Range("G20").Value = "➖ Neutre"

# Row 21: VIVO ENERGY CI (SHEC)
$ws1.Range("A21").Value = "VIVO ENERGY CI (SHEC)"
$ws1.Range("B21").Value = 1
$ws1.Range("C21").Value = 0
$ws1.Range("D21").Value = 4.39
$ws1.Range("E21").Value = 4.39
$ws1.Range("F21").Value = "🟡 Observer"
$ws1.Range("G21").Value = "➖ Neutre"

# Row 22: AFRICA GLOBAL LOGISTICS CI (SDSC)
$ws1.Range("A22").Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws1.Range("B22").Value = 2
$ws1.Range("C22").Value = 1
$ws1.Range("D22").Value = 4.28
$ws1.Range("E22").Value = -4.34
$ws1.Range("F22").Value = "🟡 Observer"
$ws1.Range("G22").Value = "👀 À surveiller"

# Row 23: CORIS BANK INTERNATIONAL (CBIBF)
$ws1.Range("A23").Value = "CORIS BANK INTERNATIONAL (CBIBF)"
$ws1.Range("B23").Value = 1
$ws1.Range("C23").Value = 0
$ws1.Range("D23").Value = 4.2
$ws1.Range("E23").Value = 4.2
$ws1.Range("F23").Value = "🟡 Observer"
$ws1.Range("G23").Value = "➖ Neutre"

# Row 24: ECOBANK TRANS. INCORP. TG (ETIT)
$ws1.Range("A24").Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Range("B24").Value = 1
$ws1.Range("C24").Value = 1
$ws1.Range("D24").Value = 0.18
$ws1.Range("E24").Value = -4.17
$ws1.Range("F24").Value = "🟡 Observer"
$ws1.Range("G24").Value = "👀 À surveiller"

# Row 25: SAPH CI (SPHC)
$ws1.Range("A25").Value = "SAPH CI (SPHC)"
$ws1.Range("B25").Value = 0
$ws1.Range("C25").Value = 1
$ws1.Range("D25").Value = -1.31
$ws1.Range("E25").Value = -1.31
$ws1.Range("F25").Value = "🟡 Observer"
$ws1.Range("G25").Value = "➖ Neutre"

# Row 26: LOTERIE NATIONALE DU BENIN (LNBB)
$ws1.Range("A26").Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws1.Range("B26").Value = 0
$ws1.Range("C26").Value = 1
$ws1.Range("D26").Value = -1.38
$ws1.Range("E26").Value = -1.38
$ws1.Range("F26").Value = "🟡 Observer"
$ws1.Range("G26").Value = "➖ Neutre"

# Row 27: SICOR CI (SICC)
$ws1.Range("A27").Value = "SICOR CI (SICC)"
$ws1.Range("B27").Value = 1
$ws1.Range("C27").Value = 1
$ws1.Range("D27").Value = -1.51
$ws1.Range("E27").Value = -7.48
$ws1.Range("F27").Value = "🟡 Observer"
$ws1.Range("G27").Value = "👀 À surveiller"

# Row 28: UNIWAX CI (UNXC)
$ws1.Range("A28").Value = "UNIWAX CI (UNXC)"
$ws1.Range("B28").Value = 1
$ws1.Range("C28").Value = 1
$ws1.Range("D28").Value = -1.73
$ws1.Range("E28").Value = 5.61
$ws1.Range("F28").Value = "🟡 Observer"
$ws1.Range("G28").Value = "👀 À surveiller"

# Row 29: BANK OF AFRICA BF (BOABF)
$ws1.Range("A29").Value = "BANK OF AFRICA BF (BOABF)"
$ws1.Range("B29").Value = 1
$ws1.Range("C29").Value = 1
$ws1.Range("D29").Value = -1.92
$ws1.Range("E29").Value = 3.9
$ws1.Range("F29").Value = "🟡 Observer"
$ws1.Range("G29").Value = "👀 À surveiller"

# Row 30: TRACTAFRIC MOTORS CI (PRSC)
$ws1.Range("A30").Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Range("B30").Value = 0
$ws1.Range("C30").Value = 1
$ws1.Range("D30").Value = -1.98
$ws1.Range("E30").Value = -1.98
$ws1.Range("F30").Value = "🟡 Observer"
$ws1.Range("G30").Value = "➖ Neutre"

# Row 31: SICABLE CI (CABC)
$ws1.Range("A31").Value = "SICABLE CI (CABC)"
$ws1.Range("B31").Value = 0
$ws1.Range("C31").Value = 1
$ws1.Range("D31").Value = -1.99
$ws1.Range("E31").Value = -1.99
$ws1.Range("F31").Value = "🟡 Observer"
$ws1.Range("G31").Value = "➖ Neutre"

# Row 32: TOTALENERGIES MARKETING CI (TTLC)
$ws1.Range("A32").Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Range("B32").Value = 0
$ws1.Range("C32").Value = 1
$ws1.Range("D32").Value = -2.08
$ws1.Range("E32").Value = -2.08
$ws1.Range("F32").Value = "🟡 Observer"
$ws1.Range("G32").Value = "➖ Neutre"

# Row 33: BERNABE CI (BNBC)
$ws1.Range("A33").Value = "BERNABE CI (BNBC)"
$ws1.Range("B33").Value = 1
$ws1.Range("C33").Value = 3
$ws1.Range("D33").Value = -2.59
$ws1.Range("E33").Value = 6.79
$ws1.Range("F33").Value = "🟡 Observer"
$ws1.Range("G33").Value = "👀 À surveiller"

# Row 34: ORAGROUP TOGO (ORGT)
$ws1.Range("A34").Value = "ORAGROUP TOGO (ORGT)"
$ws1.Range("B34").Value = 0
$ws1.Range("C34").Value = 1
$ws1.Range("D34").Value = -2.93
$ws1.Range("E34").Value = -2.93
$ws1.Range("F34").Value = "🟡 Observer"
$ws1.Range("G34").Value = "➖ Neutre"

# Row 35: NESTLE CI (NTLC)
$ws1.Range("A35").Value = "NESTLE CI (NTLC)"
$ws1.Range("B35").Value = 0
$ws1.Range("C35").Value = 1
$ws1.Range("D35").Value = -3.47
$ws1.Range("E35").Value = -3.47
$ws1.Range("F35").Value = "🟡 Observer"
$ws1.Range("G35").Value = "➖ Neutre"

# Row 36: FILTISAC CI (FTSC)
$ws1.Range("A36").Value = "FILTISAC CI (FTSC)"
$ws1.Range("B36").Value = 0
$ws1.Range("C36").Value = 2
$ws1.Range("D36").Value = -4.22
$ws1.Range("E36").Value = -2.55
$ws1.Range("F36").Value = "🟡 Observer"
$ws1.Range("G36").Value = "➖ Neutre"

# Row 37: SAFCA CI (SAFC)
$ws1.Range("A37").Value = "SAFCA CI (SAFC)"
$ws1.Range("B37").Value = 0
$ws1.Range("C37").Value = 1
$ws1.Range("D37").Value = -6.22
$ws1.Range("E37").Value = -6.22
$ws1.Range("F37").Value = "🟡 Observer"
$ws1.Range("G37").Value = "➖ Neutre"

# Row 38: SITAB CI (STBC)
$ws1.Range("A38").Value = "SITAB CI (STBC)"
$ws1.Range("B38").Value = 0
$ws1.Range("C38").Value = 1
$ws1.Range("D38").Value = -7.48
$ws1.Range("E38").Value = -7.48
$ws1.Range("F38").Value = "🟡 Observer"
$ws1.Range("G38").Value = "➖ Neutre"

# --- Sheet 2: Top_YTD ---
$ws2 = $wb.Worksheets.Item("Top_YTD")

# Row 2: BRVM - CONSOMMATION DISCRETIONNAIRE
$ws2.Range("A2").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws2.Range("B2").Value = 6014.91

# Row 3: BRVM - INDUSTRIELS
$ws2.Range("A3").Value = "BRVM - INDUSTRIELS"
$ws2.Range("B3").Value = 4361.86

# Row 4: BRVM - SERVICES FINANCIERS
$ws2.Range("A4").Value = "BRVM - SERVICES FINANCIERS"
$ws2.Range("B4").Value = 4165.77

# Row 5: BRVM-PRESTIGE
$ws2.Range("A5").Value = "BRVM-PRESTIGE"
$ws2.Range("B5").Value = 3749.77

# Row 6: BRVM - SERVICES PUBLICS
$ws2.Range("A6").Value = "BRVM - SERVICES PUBLICS"
$ws2.Range("B6").Value = 2157.67

# Row 7: BRVM - ENERGIE
$ws2.Range("A7").Value = "BRVM - ENERGIE"
$ws2.Range("B7").Value = 2152.42

# Row 8: BRVM - TELECOMMUNICATIONS
$ws2.Range("A8").Value = "BRVM - TELECOMMUNICATIONS"
$ws2.Range("B8").Value = 1380.56

# Row 9: BRVM-PRINCIPAL     (**)
$ws2.Range("A9").Value = "BRVM-PRINCIPAL     (**)"
$ws2.Range("B9").Value = 236.54

# Row 10: BRVM-PRINCIPAL
$ws2.Range("A10").Value = "BRVM-PRINCIPAL"
$ws2.Range("B10").Value = 236.03

# Row 11: BRVM - CONSOMMATION DE BASE
$ws2.Range("A11").Value = "BRVM - CONSOMMATION DE BASE"
$ws2.Range("B11").Value = 233.7
